$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format before writing, so that
# numeric-looking strings (e.g. "0.998") are stored as text, matching the
# original inlineStr cells rather than being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '49.962.82'
$ws.Range('E2').Value = '  +3.75%  '
$ws.Range('D3').Value = '2.627.73'
$ws.Range('E3').Value = '  +5.07%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '326.88'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '110.19'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '0.560'
$ws.Range('E9').Value = '  +4.12%  '
$ws.Range('D10').Value = '40.28'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('D11').Value = '20.63'
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').Value = '0.0820'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').Value = '7.28'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('D15').Value = '3.034.66'
$ws.Range('E15').Value = '  +4.84%  '
$ws.Range('D16').Value = '2.617.62'
$ws.Range('E16').Value = '  +4.54%  '
$ws.Range('D17').Value = '0.876'
$ws.Range('E17').Value = '  +4.41%  '
$ws.Range('D18').Value = '49.826.46'
$ws.Range('E18').Value = '  +3.78%  '
$ws.Range('D19').Value = '3.07'
$ws.Range('E19').Value = '  +12.21%  '
$ws.Range('D20').Value = '13.34'
$ws.Range('E20').Value = '  +1.86%  '
$ws.Range('D21').Value = '6.83'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('D23').Value = '72.75'
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').Value = '279.17'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('E25').Value = '  +1.86%  '
$ws.Range('D26').Value = '26.57'
$ws.Range('E26').Value = '  +3.39%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '36.83'
$ws.Range('E28').Value = '  +4.87%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '9.96'
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('D31').Value = '0.144'
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').Value = '49.89'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').Value = '19.73'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '5.44'
$ws.Range('E34').Value = '  +2.09%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').Value = '0.0792'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  +5.76%  '
$ws.Range('D38').Value = '4.76'
$ws.Range('E38').Value = '  +2.15%  '
$ws.Range('E39').Value = '  +7.52%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.112'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '123.22'
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('D42').Value = '22.41'
$ws.Range('E42').Value = '  +4.37%  '
$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').Value = '0.0314'
$ws.Range('E44').Value = '  +4.48%  '
$ws.Range('D45').Value = '3.36'
$ws.Range('E45').Value = '  +5.87%  '
$ws.Range('D46').Value = '2.054.60'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('D47').Value = '2.32'
$ws.Range('E47').Value = '  +16.07%  '
$ws.Range('E48').Value = '  +8.53%  '
$ws.Range('D49').Value = '9.01'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').Value = '5.36'
$ws.Range('E50').Value = '  +3.57%  '
$ws.Range('D51').Value = '81.65'
$ws.Range('E51').Value = '  +1.66%  '

# Restore the default (General) cell format now that the text values are
# committed as strings, so formatting matches the original workbook.
$ws.Range("D2:E51").ClearFormats()
